$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("MoCo")
$ws = $wb.Worksheets.Item("Res50FC")

$dst = $ws.Range("Z32")
# Step 1: border only, on blank cell -> should reuse borderId1
$dst.Borders.LineStyle = 1
$dst.Borders.Weight = 2
$dst.Borders.ColorIndex = -4105

# Step 2: paste fill but exclude borders via PasteSpecial Paste type = xlPasteAllExceptBorders (7)
$ws2.Range("A2").Copy()
$dst.PasteSpecial(7)
